# Update NATMI LR-pair TPM-derived metrics for Ccl11-Ccr2 (OldD7) per new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 0.738254
$ws.Cells.Item(2, 8).Value = 2.214762
$ws.Cells.Item(2, 9).Value = 0.005691320045803731
$ws.Cells.Item(2, 10).Value = 0.005691320045803731
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.1278363333333334
$ws.Cells.Item(2, 14).Value = 0.383509
$ws.Cells.Item(2, 15).Value = 0.002480915078704262
$ws.Cells.Item(2, 16).Value = 0.002480915078704262
$ws.Cells.Item(2, 17).Value = 0.09437568442866669
$ws.Cells.Item(2, 18).Value = 0.8493811598580001
$ws.Cells.Item(2, 19).Value = 0.00001411968171936631
$ws.Cells.Item(2, 20).Value = 0.00001411968171936631

# Row 3
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 0.738254
$ws.Cells.Item(3, 8).Value = 2.214762
$ws.Cells.Item(3, 9).Value = 0.005691320045803731
$ws.Cells.Item(3, 10).Value = 0.005691320045803731
$ws.Cells.Item(3, 15).Value = 0.002269935507489869
$ws.Cells.Item(3, 16).Value = 0.002269935507489869
$ws.Cells.Item(3, 17).Value = 0.08634987910999999
$ws.Cells.Item(3, 18).Value = 0.7771489119899999
$ws.Cells.Item(3, 19).Value = 0.00001291892945645875
$ws.Cells.Item(3, 20).Value = 0.00001291892945645875

# Row 4
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.738254
$ws.Cells.Item(4, 8).Value = 2.214762
$ws.Cells.Item(4, 9).Value = 0.005691320045803731
$ws.Cells.Item(4, 10).Value = 0.005691320045803731
$ws.Cells.Item(4, 11).Value = 2.0
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.04120633333333334
$ws.Cells.Item(4, 14).Value = 0.123619
$ws.Cells.Item(4, 15).Value = 0.0007996898146180199
$ws.Cells.Item(4, 16).Value = 0.0007996898146180199
$ws.Cells.Item(4, 17).Value = 0.03042074040866667
$ws.Cells.Item(4, 18).Value = 0.273786663678
$ws.Cells.Item(4, 19).Value = 0.000004551290672360606
$ws.Cells.Item(4, 20).Value = 0.000004551290672360606

# Row 5
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 0.738254
$ws.Cells.Item(5, 8).Value = 2.214762
$ws.Cells.Item(5, 9).Value = 0.005691320045803731
$ws.Cells.Item(5, 10).Value = 0.005691320045803731
$ws.Cells.Item(5, 13).Value = 51.241888
$ws.Cells.Item(5, 14).Value = 153.725664
$ws.Cells.Item(5, 15).Value = 0.9944494595991877
$ws.Cells.Item(5, 16).Value = 0.9944494595991878
$ws.Cells.Item(5, 17).Value = 37.82952878355199
$ws.Cells.Item(5, 18).Value = 340.465759051968
$ws.Cells.Item(5, 19).Value = 0.005659730143955545
$ws.Cells.Item(5, 20).Value = 0.005659730143955545

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9440493064670392
$ws.Cells.Item(6, 10).Value = 0.9440493064670391
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 0.1278363333333334
$ws.Cells.Item(6, 14).Value = 0.383509
$ws.Cells.Item(6, 15).Value = 0.002480915078704262
$ws.Cells.Item(6, 16).Value = 0.002480915078704262
$ws.Cells.Item(6, 17).Value = 15.654593084767
$ws.Cells.Item(6, 18).Value = 140.891337762903
$ws.Cells.Item(6, 19).Value = 0.002342106159454379
$ws.Cells.Item(6, 20).Value = 0.002342106159454379

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9440493064670392
$ws.Cells.Item(7, 10).Value = 0.9440493064670391
$ws.Cells.Item(7, 15).Value = 0.002269935507489869
$ws.Cells.Item(7, 16).Value = 0.002269935507489869
$ws.Cells.Item(7, 19).Value = 0.002142931041570717
$ws.Cells.Item(7, 20).Value = 0.002142931041570717

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9440493064670392
$ws.Cells.Item(8, 10).Value = 0.9440493064670391
$ws.Cells.Item(8, 11).Value = 2.0
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.04120633333333334
$ws.Cells.Item(8, 14).Value = 0.123619
$ws.Cells.Item(8, 15).Value = 0.0007996898146180199
$ws.Cells.Item(8, 16).Value = 0.0007996898146180199
$ws.Cells.Item(8, 17).Value = 5.046048834697001
$ws.Cells.Item(8, 18).Value = 45.41443951227301
$ws.Cells.Item(8, 19).Value = 0.0007549466148788969
$ws.Cells.Item(8, 20).Value = 0.0007549466148788968

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9440493064670392
$ws.Cells.Item(9, 10).Value = 0.9440493064670391
$ws.Cells.Item(9, 13).Value = 51.241888
$ws.Cells.Item(9, 14).Value = 153.725664
$ws.Cells.Item(9, 15).Value = 0.9944494595991877
$ws.Cells.Item(9, 16).Value = 0.9944494595991878
$ws.Cells.Item(9, 17).Value = 6274.983681232033
$ws.Cells.Item(9, 18).Value = 56474.85313108829
$ws.Cells.Item(9, 19).Value = 0.9388093226511351
$ws.Cells.Item(9, 20).Value = 0.9388093226511351

# Row 10
$ws.Cells.Item(10, 7).Value = 5.698467
$ws.Cells.Item(10, 8).Value = 17.095401
$ws.Cells.Item(10, 9).Value = 0.0439304080539368
$ws.Cells.Item(10, 10).Value = 0.04393040805393679
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.1278363333333334
$ws.Cells.Item(10, 14).Value = 0.383509
$ws.Cells.Item(10, 15).Value = 0.002480915078704262
$ws.Cells.Item(10, 16).Value = 0.002480915078704262
$ws.Cells.Item(10, 17).Value = 0.7284711269010001
$ws.Cells.Item(10, 18).Value = 6.556240142109
$ws.Cells.Item(10, 19).Value = 0.000108987611754643
$ws.Cells.Item(10, 20).Value = 0.000108987611754643

# Row 11
$ws.Cells.Item(11, 7).Value = 5.698467
$ws.Cells.Item(11, 8).Value = 17.095401
$ws.Cells.Item(11, 9).Value = 0.0439304080539368
$ws.Cells.Item(11, 10).Value = 0.04393040805393679
$ws.Cells.Item(11, 15).Value = 0.002269935507489869
$ws.Cells.Item(11, 16).Value = 0.002269935507489869
$ws.Cells.Item(11, 17).Value = 0.6665211926549999
$ws.Cells.Item(11, 18).Value = 5.998690733894999
$ws.Cells.Item(11, 19).Value = 0.00009971919310015003
$ws.Cells.Item(11, 20).Value = 0.00009971919310015002

# Row 12
$ws.Cells.Item(12, 7).Value = 5.698467
$ws.Cells.Item(12, 8).Value = 17.095401
$ws.Cells.Item(12, 9).Value = 0.0439304080539368
$ws.Cells.Item(12, 10).Value = 0.04393040805393679
$ws.Cells.Item(12, 11).Value = 2.0
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.04120633333333334
$ws.Cells.Item(12, 14).Value = 0.123619
$ws.Cells.Item(12, 15).Value = 0.0007996898146180199
$ws.Cells.Item(12, 16).Value = 0.0007996898146180199
$ws.Cells.Item(12, 17).Value = 0.234812930691
$ws.Cells.Item(12, 18).Value = 2.113316376219
$ws.Cells.Item(12, 19).Value = 0.00003513069987274668
$ws.Cells.Item(12, 20).Value = 0.00003513069987274668

# Row 13
$ws.Cells.Item(13, 7).Value = 5.698467
$ws.Cells.Item(13, 8).Value = 17.095401
$ws.Cells.Item(13, 9).Value = 0.0439304080539368
$ws.Cells.Item(13, 10).Value = 0.04393040805393679
$ws.Cells.Item(13, 13).Value = 51.241888
$ws.Cells.Item(13, 14).Value = 153.725664
$ws.Cells.Item(13, 15).Value = 0.9944494595991877
$ws.Cells.Item(13, 16).Value = 0.9944494595991878
$ws.Cells.Item(13, 17).Value = 292.000207785696
$ws.Cells.Item(13, 18).Value = 2628.001870071264
$ws.Cells.Item(13, 19).Value = 0.04368657054920926
$ws.Cells.Item(13, 20).Value = 0.04368657054920925

# Row 14
$ws.Cells.Item(14, 7).Value = 0.8209666666666666
$ws.Cells.Item(14, 8).Value = 2.4629
$ws.Cells.Item(14, 9).Value = 0.006328965433220369
$ws.Cells.Item(14, 10).Value = 0.006328965433220369
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 0.1278363333333334
$ws.Cells.Item(14, 14).Value = 0.383509
$ws.Cells.Item(14, 15).Value = 0.002480915078704262
$ws.Cells.Item(14, 16).Value = 0.002480915078704262
$ws.Cells.Item(14, 17).Value = 0.1049493684555556
$ws.Cells.Item(14, 18).Value = 0.9445443161
$ws.Cells.Item(14, 19).Value = 0.00001570162577587447
$ws.Cells.Item(14, 20).Value = 0.00001570162577587447

# Row 15
$ws.Cells.Item(15, 7).Value = 0.8209666666666666
$ws.Cells.Item(15, 8).Value = 2.4629
$ws.Cells.Item(15, 9).Value = 0.006328965433220369
$ws.Cells.Item(15, 10).Value = 0.006328965433220369
$ws.Cells.Item(15, 15).Value = 0.002269935507489869
$ws.Cells.Item(15, 16).Value = 0.002269935507489869
$ws.Cells.Item(15, 17).Value = 0.09602436616666665
$ws.Cells.Item(15, 18).Value = 0.8642192954999999
$ws.Cells.Item(15, 19).Value = 0.00001436634336254291
$ws.Cells.Item(15, 20).Value = 0.00001436634336254291

# Row 16
$ws.Cells.Item(16, 7).Value = 0.8209666666666666
$ws.Cells.Item(16, 8).Value = 2.4629
$ws.Cells.Item(16, 9).Value = 0.006328965433220369
$ws.Cells.Item(16, 10).Value = 0.006328965433220369
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.04120633333333334
$ws.Cells.Item(16, 14).Value = 0.123619
$ws.Cells.Item(16, 15).Value = 0.0007996898146180199
$ws.Cells.Item(16, 16).Value = 0.0007996898146180199
$ws.Cells.Item(16, 17).Value = 0.03382902612222222
$ws.Cells.Item(16, 18).Value = 0.3044612351
$ws.Cells.Item(16, 19).Value = 0.000005061209194015853
$ws.Cells.Item(16, 20).Value = 0.000005061209194015853

# Row 17
$ws.Cells.Item(17, 7).Value = 0.8209666666666666
$ws.Cells.Item(17, 8).Value = 2.4629
$ws.Cells.Item(17, 9).Value = 0.006328965433220369
$ws.Cells.Item(17, 10).Value = 0.006328965433220369
$ws.Cells.Item(17, 13).Value = 51.241888
$ws.Cells.Item(17, 14).Value = 153.725664
$ws.Cells.Item(17, 15).Value = 0.9944494595991877
$ws.Cells.Item(17, 16).Value = 0.9944494595991878
$ws.Cells.Item(17, 17).Value = 42.06788198506666
$ws.Cells.Item(17, 18).Value = 378.6109378656
$ws.Cells.Item(17, 19).Value = 0.006293836254887935
$ws.Cells.Item(17, 20).Value = 0.006293836254887936
